$x = 1 + 2
Write-Host $x
$arr = @(1,2,3)
Write-Host $arr.Count
if ($x -eq 3) { Write-Host "yes" }
for ($i=0; $i -lt 3; $i++) { Write-Host $i }
function Foo($a) { return $a * 2 }
Write-Host (Foo 5)
